# Auto-generated Excel COM-interop script applying a scheduled market-data
# refresh to the Hades_Profits leve-profit workbook. Each block updates the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for a single
# leve row on a single crafting-job sheet, using the values captured in the
# upstream commit "chore: update Sheets via scheduled runner".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1696.25
$ws.Range("J17").Value = 1696.25
$ws.Range("L17").Value = 5088.75
$ws.Range("N17").Value = -5424.75
# Row 129
$ws.Range("H129").Value = 788.22
$ws.Range("I129").Value = 559.1818
$ws.Range("J129").Value = 968.1786
$ws.Range("K129").Value = 1677.5454
$ws.Range("L129").Value = 2904.5358
$ws.Range("M129").Value = 3322.4546
$ws.Range("N129").Value = -12904.5358
# Row 134
$ws.Range("H134").Value = 37846.152
$ws.Range("J134").Value = 37846.152
$ws.Range("L134").Value = 37846.152
$ws.Range("N134").Value = -47986.152
# Row 138
$ws.Range("H138").Value = 4258196
$ws.Range("I138").Value = 2948
$ws.Range("J138").Value = 4880915.5
$ws.Range("K138").Value = 8844
$ws.Range("L138").Value = 14642746.5
$ws.Range("M138").Value = -3704
$ws.Range("N138").Value = -14653026.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1438.6666
$ws.Range("I45").Value = 1349.7142
$ws.Range("J45").Value = 1750
$ws.Range("K45").Value = 1349.7142
$ws.Range("L45").Value = 1750
$ws.Range("M45").Value = -972.7141999999999
$ws.Range("N45").Value = -2504
# Row 61
$ws.Range("H61").Value = 167001310
$ws.Range("I61").Value = 200201580
$ws.Range("J61").Value = 1000000
$ws.Range("K61").Value = 200201580
$ws.Range("L61").Value = 1000000
$ws.Range("M61").Value = -200201368
$ws.Range("N61").Value = -1000424
# Row 74
$ws.Range("H74").Value = 6212078.5
$ws.Range("I74").Value = 9834911
$ws.Range("J74").Value = 53262.4
$ws.Range("K74").Value = 9834911
$ws.Range("L74").Value = 53262.4
$ws.Range("M74").Value = -9834037
$ws.Range("N74").Value = -55010.4
# Row 77
$ws.Range("H77").Value = 6212078.5
$ws.Range("I77").Value = 9834911
$ws.Range("J77").Value = 53262.4
$ws.Range("K77").Value = 49174555
$ws.Range("L77").Value = 266312
$ws.Range("M77").Value = -49170187
$ws.Range("N77").Value = -275048
# Row 110
$ws.Range("H110").Value = 1676.7142
$ws.Range("I110").Value = 1111.625
$ws.Range("K110").Value = 1111.625
$ws.Range("M110").Value = 933.375
# Row 136
$ws.Range("H136").Value = 167001310
$ws.Range("I136").Value = 200201580
$ws.Range("J136").Value = 1000000
$ws.Range("K136").Value = 600604740
$ws.Range("L136").Value = 3000000
$ws.Range("M136").Value = -600602190
$ws.Range("N136").Value = -3005100

$ws = $wb.Worksheets.Item("CRP")
# Row 52
$ws.Range("H52").Value = 49900
$ws.Range("J52").Value = 49900
$ws.Range("L52").Value = 49900
$ws.Range("N52").Value = -50488
# Row 137
$ws.Range("H137").Value = 41999.8
$ws.Range("J137").Value = 41999.8
$ws.Range("L137").Value = 41999.8
$ws.Range("N137").Value = -52199.8
# Row 141
$ws.Range("H141").Value = 47364.285
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 47364.285
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 47364.285
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -57724.285

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 613.84375
$ws.Range("I113").Value = 559.6087
$ws.Range("J113").Value = 752.44446
$ws.Range("K113").Value = 1678.8261
$ws.Range("L113").Value = 2257.33338
$ws.Range("M113").Value = 491.1739
$ws.Range("N113").Value = -6597.33338
# Row 121
$ws.Range("H121").Value = 76837910
$ws.Range("J121").Value = 131721870
$ws.Range("L121").Value = 395165610
$ws.Range("N121").Value = -395168230
# Row 131
$ws.Range("H131").Value = 1011.0833
$ws.Range("J131").Value = 1181.4445
$ws.Range("L131").Value = 3544.3335
$ws.Range("N131").Value = -13624.3335
# Row 134
$ws.Range("H134").Value = 4332.7827
$ws.Range("I134").Value = 2815.875
$ws.Range("J134").Value = 7800
$ws.Range("K134").Value = 8447.625
$ws.Range("L134").Value = 23400
$ws.Range("M134").Value = -3377.625
$ws.Range("N134").Value = -33540

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1000
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -1224
# Row 22
$ws.Range("H22").Value = 673
$ws.Range("I22").Value = 563.6429000000001
$ws.Range("J22").Value = 782.3570999999999
$ws.Range("K22").Value = 563.6429000000001
$ws.Range("L22").Value = 782.3570999999999
$ws.Range("M22").Value = -268.6429000000001
$ws.Range("N22").Value = -1372.3571
# Row 27
$ws.Range("H27").Value = 673
$ws.Range("I27").Value = 563.6429000000001
$ws.Range("J27").Value = 782.3570999999999
$ws.Range("K27").Value = 563.6429000000001
$ws.Range("L27").Value = 782.3570999999999
$ws.Range("M27").Value = -456.6429000000001
$ws.Range("N27").Value = -996.3570999999999
# Row 61
$ws.Range("H61").Value = 2568.12
$ws.Range("I61").Value = 2655.4443
$ws.Range("J61").Value = 2343.5715
$ws.Range("K61").Value = 2655.4443
$ws.Range("L61").Value = 2343.5715
$ws.Range("M61").Value = -2453.4443
$ws.Range("N61").Value = -2747.5715
# Row 82
$ws.Range("H82").Value = 2802.8
$ws.Range("I82").Value = 2002
$ws.Range("K82").Value = 2002
$ws.Range("M82").Value = -1641
# Row 85
$ws.Range("H85").Value = 2802.8
$ws.Range("I85").Value = 2002
$ws.Range("K85").Value = 2002
$ws.Range("M85").Value = -754
# Row 93
$ws.Range("H93").Value = 1209.8334
$ws.Range("I93").Value = 1207.5555
$ws.Range("K93").Value = 1207.5555
$ws.Range("M93").Value = 40.44450000000006
# Row 100
$ws.Range("H100").Value = 1816.8334
$ws.Range("I100").Value = 1620.6
$ws.Range("J100").Value = 1892.3077
$ws.Range("K100").Value = 1620.6
$ws.Range("L100").Value = 1892.3077
$ws.Range("M100").Value = -1079.6
$ws.Range("N100").Value = -2974.3077
# Row 113
$ws.Range("H113").Value = 2568.12
$ws.Range("I113").Value = 2655.4443
$ws.Range("J113").Value = 2343.5715
$ws.Range("K113").Value = 2655.4443
$ws.Range("L113").Value = 2343.5715
$ws.Range("M113").Value = -485.4443000000001
$ws.Range("N113").Value = -6683.5715
# Row 132
$ws.Range("H132").Value = 108600.2
$ws.Range("I132").Value = 7200.6
$ws.Range("J132").Value = 209999.8
$ws.Range("K132").Value = 21601.8
$ws.Range("L132").Value = 629999.3999999999
$ws.Range("M132").Value = -19071.8
$ws.Range("N132").Value = -635059.3999999999
# Row 136
$ws.Range("H136").Value = 36955.117
$ws.Range("I136").Value = 23550.348
$ws.Range("J136").Value = 102638.5
$ws.Range("K136").Value = 70651.04400000001
$ws.Range("L136").Value = 307915.5
$ws.Range("M136").Value = -68101.04400000001
$ws.Range("N136").Value = -313015.5

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 427.625
$ws.Range("I107").Value = 401.07693
$ws.Range("J107").Value = 459
$ws.Range("K107").Value = 1203.23079
$ws.Range("L107").Value = 1377
$ws.Range("M107").Value = 716.7692099999999
$ws.Range("N107").Value = -5217
# Row 113
$ws.Range("H113").Value = 726.5517
$ws.Range("I113").Value = 891
$ws.Range("J113").Value = 493.58334
$ws.Range("K113").Value = 2673
$ws.Range("L113").Value = 1480.75002
$ws.Range("M113").Value = -503
$ws.Range("N113").Value = -5820.750019999999

